$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.682.67'
$ws.Range('D3').Value = '2.577.07'
$ws.Range('E3').Value = '  +1.19%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '581.14'
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '145.16'
$ws.Range('E6').Value = '  -1.12%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.593'
$ws.Range('E8').Value = '  +1.64%  '
$ws.Range('E9').Value = '  +0.89%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '5.59'
$ws.Range('E10').Value = '  +0.64%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.151'
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('E12').Value = '  -0.58%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '27.04'
$ws.Range('E13').Value = '  -1.61%  '
$ws.Range('D14').Value = '3.039.91'
$ws.Range('E14').Value = '  +1.12%  '
$ws.Range('D15').Value = '62.584.31'
$ws.Range('E15').Value = '  -0.55%  '
$ws.Range('E16').Value = '  +1.51%  '
$ws.Range('D17').Value = '2.583.54'
$ws.Range('E17').Value = '  +1.57%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.21'
$ws.Range('E18').Value = '  -1.14%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '338.89'
$ws.Range('E19').Value = '  +0.39%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.36'
$ws.Range('E20').Value = '  +0.79%  '
$ws.Range('E21').Value = '  -1.34%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '67.06'
$ws.Range('E23').Value = '  +2.15%  '
$ws.Range('D24').Value = '2.703.45'
$ws.Range('E24').Value = '  +0.93%  '
$ws.Range('E25').Value = '  -2.40%  '
$ws.Range('E26').Value = '  -1.29%  '
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('B28').Value = 'SuiNetwork'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.48'
$ws.Range('E28').Value = '  -0.33%  '
$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.91'
$ws.Range('E29').Value = '  +2.87%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.25'
$ws.Range('E30').Value = '  -1.18%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.92'
$ws.Range('E31').Value = '  -3.25%  '
$ws.Range('D32').Value = '0.0₃0814'
$ws.Range('E32').Value = '  +0.09%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '461.10'
$ws.Range('E33').Value = '  +10.57%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '177.03'
$ws.Range('E34').Value = '  -0.56%  '
$ws.Range('E35').Value = '  +3.37%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.401'
$ws.Range('E37').Value = '  +0.21%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '18.91'
$ws.Range('E38').Value = '  -1.06%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.50'
$ws.Range('E39').Value = '  +3.10%  '
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.69'
$ws.Range('E41').Value = '  -3.17%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '158.49'
$ws.Range('E42').Value = '  +5.22%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.72'
$ws.Range('E43').Value = '  -1.33%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '21.00'
$ws.Range('E44').Value = '  +1.24%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.627'
$ws.Range('E45').Value = '  +4.20%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0536'
$ws.Range('E46').Value = '  -0.69%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0964'
$ws.Range('E47').Value = '  -0.51%  '
$ws.Range('E48').Value = '  -1.29%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '18.15'
$ws.Range('E49').Value = '  -0.67%  '
$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '11.42'
$ws.Range('E50').Value = '  +1.07%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.71'
$ws.Range('E51').Value = '  -0.02%  '
